# Rede de serviços.xlsx - update equipment names ("Casa da Mulher" rebrand)
# and correct the location of "Casa da Mulher Cachoeirinha".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Rename facilities (column D = nm_equipam) -----------------------------
# Old "Centro de Defesa e de Convivência da Mulher / Centro de Cidadania da
# Mulher / Centro de Referência e Cidadania da Mulher - ..." names are being
# simplified to the "Casa da Mulher ..." naming convention.

$ws.Range("D4").Value  = "Casa da Mulher Mulheres Vivas"
$ws.Range("D5").Value  = "Casa da Mulher Cidinha Kopcak"
$ws.Range("D6").Value  = "Casa da Mulher Anastácia"
$ws.Range("D7").Value  = "Casa da Mulher Marielle Franco"
$ws.Range("D8").Value  = "Casa da Mulher Marcia Martins"
$ws.Range("D9").Value  = "Casa da Mulher Viviane dos Santos"
$ws.Range("D10").Value = "Casa da Mulher Sônia Maria Batista"
$ws.Range("D11").Value = "Casa da Mulher Crê-ser"
$ws.Range("D17").Value = "Casa da Mulher Isabel Projeto Naná Serafim"
$ws.Range("D18").Value = "Casa da Mulher Sofia"
$ws.Range("D19").Value = "Casa da Mulher Centro de Integração Social da Mulher"
$ws.Range("D20").Value = "Casa da Mulher Margarida Maria Alves"
$ws.Range("D23").Value = "Casa da Mulher Capela do Socorro"
$ws.Range("D24").Value = "Casa da Mulher Itaquera"
$ws.Range("D25").Value = "Casa da Mulher Parelheiros"
$ws.Range("D26").Value = "Casa da Mulher Perus"
$ws.Range("D27").Value = "Casa da Mulher  Santo Amaro"
$ws.Range("D32").Value = "Casa da Mulher Mariás"
$ws.Range("D33").Value = "Casa da Mulher Zizi"
$ws.Range("D34").Value = "Casa da Mulher Espaço Francisca Franco"
$ws.Range("D45").Value = "Casa da Mulher 25 de Março"
$ws.Range("D46").Value = "Casa da Mulher  Brasilândia"
$ws.Range("D47").Value = "Casa da Mulher Cachoeirinha"
$ws.Range("D48").Value = "Casa da Mulher Săo Miguel"
$ws.Range("D49").Value = "Casa da Mulher Eliane de Grammont"
$ws.Range("D50").Value = "Casa da Mulher Maria de Lourdes Rodrigues"

# --- Correct the location of "Casa da Mulher Cachoeirinha" (row 47, -------
# cd_identif 270048): address, CEP and lat/long were all wrong. -------------
$ws.Range("K47").Value = "Avenida Deputado Emílio Carlos, 3460"
$ws.Range("L47").Value = "02721-200"
$ws.Range("M47").Value = -23.477804053634369
$ws.Range("N47").Value = -46.671271274611691
